# Edit: (1) change the table style used on the SOURCES OF FINANCE table
#       (slide 6) from the deck's custom "Table_0" style to the built-in
#       PowerPoint table style "Medium Style 2 - Accent 1", and
#       (2) swap the colour palette of the presentation's theme so the
#       slide master now uses the default Office-theme colours instead
#       of the custom "Integral" colours.

function HexToRgbInt($hex) {
    # VBA/PowerPoint RGB() packs r + g*256 + b*65536 into the long that
    # the ColorFormat/ThemeColor .RGB property expects.
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 6 -------------------------------------------
$slide6 = $p.Slides.Item(6)
for ($i = 1; $i -le $slide6.Shapes.Count; $i++) {
    $shp = $slide6.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{8F0CABCE-A337-4A32-AE70-88EA511A7471}")
    }
}

# --- 2. Theme colours ------------------------------------------------------
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink -- in that ThemeColorScheme
# index order -- switched from the "Integral" theme to the standard
# "Office Theme" palette.
$officeThemeColors = @(
    "000000", # dk1
    "FFFFFF", # lt1
    "44546A", # dk2
    "E7E6E6", # lt2
    "5B9BD5", # accent1
    "ED7D31", # accent2
    "A5A5A5", # accent3
    "FFC000", # accent4
    "4472C4", # accent5
    "70AD47", # accent6
    "0563C1", # hlink
    "954F72"  # folHlink
)

$slide1 = $p.Slides.Item(1)
$themeColors = $slide1.ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $themeColors.Item($i).RGB = HexToRgbInt($officeThemeColors[$i - 1])
}
